$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'If you are in need of an hones tGarage look no further.'
$ws.Range("F2").Value = 'Really impressed!'
$ws.Range("H2").Value = 'Not what I asked for'
$ws.Range("A3").Value = 'Twitter'
$ws.Range("B3").Value = 'my favorite website they are proven to pay me now I just work at home in this time of pandemic now they really pay thank you'
$ws.Range("C3").Value = 'AMAZING'
$ws.Range("D3").Value = '4.0'
$ws.Range("E3").Value = 'Outstanding application to feel free from threats. Enjoy to use twitter more than before with this friendly app...'
$ws.Range("F3").Value = 'comfertbility'
$ws.Range("G3").Value = 'Very good service for the affordable price! Their IPTV activation is very fast and customer support also best. It’s the best streaming service with full HD quality. My experience is great on bestiptvfrance. Highly recommended service!!!'
$ws.Range("H3").Value = 'Very good service for the affordable…'
$ws.Range("A4").Value = 'Facebook'
$ws.Range("B4").Value = 'So I joined facebook about 2 months ago and literally the only thing I did was join my town''s community page. I didn''t post anything I didn''t add anyone I didn''t write any comments and in fact I haven''t even gone on the platform in quite some time. Despite all this my account has been disabled FOUR TIMES for supposedly "violating the community standards". And somehow I''ve been on instagram for nearly a year and NEVER ONCE had this problem. Honestly as soon as it gets re-enabled again I''m just deleting my account because I''m sick of dealing with this'
$ws.Range("C4").Value = 'artificial picture of lives'
$ws.Range("E4").Value = 'THANK YOU VERY VERY MUCH FOR BEING IMPECCABLY FLAWLESS AND IMMACULATE TO THE VERY HIGHEST DEGREE!'
$ws.Range("F4").Value = 'ABSOLUTE PERFECTION GUARANTEED'
$ws.Range("G4").Value = 'ABSOLUTE PERFECTION GUARANTEED! It successfully taught me the importance of responsible communication! I became a lot more diplomatic tactful courteous and respectful! '
$ws.Range("H4").Value = 'Worst service I''ve ever gotten from a website'
$ws.Range("A5").Value = 'etyres'
$ws.Range("B5").Value = 'Great service ordered my tyres paid extra for mobile tyre fitter. Tyres came with fitter and the guy fitted them within 20 minutes.'
$ws.Range("C5").Value = 'Fantastic tyre service'
$ws.Range("D5").Value = '3.5'
$ws.Range("E5").Value = 'Abysmal!!!'
$ws.Range("F5").Value = 'Abysmal!!!'
$ws.Range("G5").Value = 'Brilliant easy to navigate website and cheapest tyres on the market will definitely purchase tyres for my other vehicles from here! :)'
$ws.Range("H5").Value = 'Happy customer :)'
$ws.Range("A6").Value = 'Garolla'
$ws.Range("B6").Value = 'The team that was sent where fantastic on time no mess explained how the door worked before they left house. Professional job all round.'
$ws.Range("C6").Value = 'The team that was sent where fantastic'
$ws.Range("D6").Value = '4.8'
$ws.Range("E6").Value = 'Very good work well pleased.'
$ws.Range("F6").Value = 'Roller doors'
$ws.Range("G6").Value = 'Delighted with Garolla Electric Roller Garage Door.  Installed on the day advised at original survey and left in a tidy and clean state following installation'
$ws.Range("H6").Value = 'Delighted with Garolla Electric Roller…'
$ws.Range("A7").Value = 'IN''n''OUT Autocentres'
$ws.Range("B7").Value = 'MOT and interim service - passed the MOT with no issues but they asked for £100 for air filter and £100 for fuel filter ( I thought filters were part of the service) they also said I need a new battery £200 - I sourced exact same battery for £97 air filter (Bosch) £18 - and just ordered fuel filter £18 and will fit them myself saving almost £270'
$ws.Range("C7").Value = 'MOT and interim service'
$ws.Range("D7").Value = '4.5'
$ws.Range("E7").Value = 'Took our Skoda Yeti on 21st February to have 2 new front brake discs and pads fitted. Charged a good price and it was ready to collect 20 minutes before stated.'
$ws.Range("F7").Value = 'More like In''n''Wait''n''Wait''n''Wait'
$ws.Range("G7").Value = 'Excellent. Excellent. Excellent! The Newark on Trent centre never fails to impress with their efficiency professionalism customer care and prices.'
$ws.Range("H7").Value = 'Excellent'
$ws.Range("A8").Value = 'Ring'
$ws.Range("B8").Value = 'Ring doorbell.lost my wifi and cant connect.tryed everything that there support suggests and now is a paper waight.i am not impressed.left meeages.tried phone contacts that not reconisable when u phone.not a good service'
$ws.Range("C8").Value = 'Ring doorbell.lost my wifi and cant…'
$ws.Range("D8").Value = '1.3'
$ws.Range("E8").Value = 'I have been waiting for a new door sensor for weeks they keep telling me they will let me know when it has been sent out but they dont send one out.'
$ws.Range("F8").Value = 'I think ring camera is great for €3 a…'
$ws.Range("G8").Value = 'I bought this system and it has been nothing but problems. I definatly dont recommend it so many faults poor customer service. '
$ws.Range("H8").Value = 'Still not replaced new faulty doorbell sensor after 2 MONTHS of waiting! Not good enough'
$ws.Range("A9").Value = 'Centralautopoint'
$ws.Range("B9").Value = 'If you are in need of an honest reliable and friendly garage you found it.'
$ws.Range("C9").Value = 'If you are in need of an hones tGarage look no further.'
$ws.Range("D9").Value = '3.6'
$ws.Range("E9").Value = 'I was not disappointed.'
$ws.Range("F9").Value = 'Really impressed!'
$ws.Range("G9").Value = 'Booked my car in for the second year for a major service and MOT this time.'
$ws.Range("H9").Value = 'Not what I asked for'
$ws.Range("A10").Value = 'E.ON SE'
$ws.Range("B10").Value = 'The wait on phone is really bad 1 half hours disgraceful.'
$ws.Range("C10").Value = 'The wait on phone is really bad 1 half…'
$ws.Range("D10").Value = '1.2'
$ws.Range("E10").Value = 'Looks like the only was is to cancel my direct debit to get them to contact me they will try and charge me for it but I''ve done everything possible so good luck with that.'
$ws.Range("F10").Value = 'Trying to find a non existant electric…'
$ws.Range("G10").Value = 'Can''t get hold of eon customer service phoned them and was on hold for an hour with no answer went on the online chat and was told to phone them tried to send them an email and it was returned went on FB messenger and was told to wait.'
$ws.Range("H10").Value = 'Can''t get hold of eon customer service'
$ws.Range("A11").Value = 'Mr Clutch'
$ws.Range("B11").Value = 'Fodi in Leyton branch best man very helpfulservice was great highly recommend !'
$ws.Range("C11").Value = 'Fodi in Leyton branch best man'
$ws.Range("D11").Value = 'Warning! We’ve detected misuse on this page'
$ws.Range("E11").Value = 'I booked my car in at the Chelmsford branch and Chris was my first contact a delightful man with perfect manners and so helpful I limped my car into the Chelmsford branch and as he promised Chris and his team had my car running perfectly by mid day I am very satisfied and will be using them for any job that comes up in the future. '
$ws.Range("F11").Value = 'I booked my car in at the Chelmsford…'
$ws.Range("G11").Value = 'Excellent service from the manager Steve I really cannot recommend them enough. Got my MOT done at the Chingford branch! Well done guys i appreciate it!'
$ws.Range("H11").Value = 'Excellent service from the manager…'
$ws.Range("A12").Value = 'The Garage Door Centre'
$ws.Range("B12").Value = 'Very pleased with the Horman roller garage door fitted today. Very professional service from the Garage Door Centre and happy to recommend. Big thank you to Scot and Justin for the fitting today.'
$ws.Range("C12").Value = 'Very pleased with the Horman roller…'
$ws.Range("D12").Value = '4.5'
$ws.Range("E12").Value = 'I was very pleased with the quality of the Hormann Georgian door for the price given. The two guys who fitted it turned up on time were friendly and professional and explained afterwards how to operate and look after the door. My old dented rusting door was removed as well at no extra cost. Compared to some local fitters'' prices I was very satisfied. Payment is required in advance; a fifty percent deposit and then the balance 3 days before fitting. This is to stop those despicable people out there who have work done and then do not pay afterwards. I was happy to pay this way and I would highly recommend this company.'
$ws.Range("F12").Value = 'Quality and value and professional'
$ws.Range("G12").Value = 'Very pleased with my new Hormann sectional garage door from The Garage Door Centre and the whole process from showroom to fitting was straightforward. The installation was completed on schedule and their team were well organised and courteous. I would recommend them.'
$ws.Range("H12").Value = 'Very pleased with my new garage door'
$ws.Range("A13").Value = 'GE'
$ws.Range("B13").Value = 'Forget about it. YOu call and can never get a person; they shunt you to the website which is equally useless.  This was to see about fixing our microwave. They clearly do not want to help you.'
$ws.Range("C13").Value = 'Forget about it'
$ws.Range("D13").Value = '1.6'
$ws.Range("E13").Value = 'We had ongoing problems since the washer was 6 months old. GE refused to acknowledge despite2-3 service calls because it wouldn''t spin or drain properly. The washer''s belt and sensor also broke at 2 1/2 years requiring service again. 1 month after our 3year warrant expired through Best Buy it stopped working completely starting with a non spin or drain again. Neither GE or Best Buy extended warranty would cover. GE customer service reps were rude and just do not care that a 3year old washer broke completely. I took this concern to the head supervisor in which she only offered to pay 50% of parts. Absolutely ridiculous. This washer was a complete lemon from the start and GE refused to acknowledge this or the real problem. My husband and I will tell everyone about this experience and never buy a GE product again!'
$ws.Range("F13").Value = 'Awful Customer Service'
$ws.Range("G13").Value = 'Do not buy the top loader-GTW465ASNWWSKU:6325628 GE'
$ws.Range("H13").Value = 'Do not buy the top…'
$ws.Range("A14").Value = 'Sheds.co.uk'
$ws.Range("B14").Value = 'Great service shed arrived exactly at the right time and was really easy to build. I also ordered some shelves and a Solar PIR light with my order. The light unfortunately was faulty but after emailing the sheds team Jenny was quick to rectify this and sent out a new PIR light replacement straight away.'
$ws.Range("C14").Value = 'Great service'
$ws.Range("D14").Value = '4.1'
$ws.Range("E14").Value = 'I would strongly recommend NOT using this company. Took money on day of order and then took 7 months to deliver the log cabin. AND THEN even managed to deliver on the wrong date! There were faults in the log cabin with no reasonable solution; and then they refused to communicate with me. Eventually I had to get my money back from Barclaycard Visa (who were excellent). Sheds.co.uk - an appalling company to do business with.'
$ws.Range("F14").Value = 'WARNING - DO NOT USE THIS COMPANY'
$ws.Range("G14").Value = 'Very helpful responses to questions asked regarding chosen shed. Thanks Ellis. Joy O'
$ws.Range("H14").Value = 'Very helpful responses to questions…'
$ws.Range("A15").Value = 'Ring'
$ws.Range("B15").Value = 'Unlike some other brands such as Canary there is no auto switching of modes ie from home to away when leaving the house. Although it does have geofence. There is no night mode either for when your asleep. You have to manually set it to away mode. The system needs to be updated online with other brands.'
$ws.Range("C15").Value = 'Unlike some other brands such as Canary…'
$ws.Range("D15").Value = '1.6'
$ws.Range("E15").Value = 'Kept losing the signal when I was setting it up. I was only two yards away from the router so what will it be like when it’s further away and outside? I use my laptop at a much further distance and have no problem with that. The app on the phone and the Fire tablet was very slow to respond sometimes not working at all even when there was a WiFi signal. On the plus side when it was working the picture quality was pretty good but that is only the good thing I can say. Returned mine for a refund.'
$ws.Range("F15").Value = 'Ring Video Doorbell'
$ws.Range("G15").Value = 'Dont buy these products i bought a outside camera it has intermitting notifications i checked on line to contact the support the online bot gave me a number i called the number it then cost me £30 because this was not a uk number. After 50 minutes on the line they have now told me i need one of their WIFI extenders because the signal is week. I was told it could be a wall blocking the signal the camera is a outdoor one so there will be a wall. I would buy another make of camera these are terrible.'
$ws.Range("H15").Value = 'Dont buy these out door camera'
$ws.Range("A16").Value = 'National Tyres and Autocare'
$ws.Range("B16").Value = 'I buy new tyres from National Tyres in Portsmouth and they fixed wrong inside place to outside ridiculous'
$ws.Range("C16").Value = 'simple and effective'
$ws.Range("D16").Value = '4.6'
$ws.Range("E16").Value = 'As always Top quality Service from National Tyres and Autocare at 317 Purley Way. '
$ws.Range("F16").Value = 'Recent Purchase'
$ws.Range("G16").Value = 'The site was very slow to load and I found the different tyre options a bit confusing.  Unable to find out if bluelight discount can be applied.'
$ws.Range("H16").Value = 'They is irresponsible'
$ws.Range("A17").Value = 'Rated People'
$ws.Range("B17").Value = 'Very disappointed '
$ws.Range("C17").Value = 'Very disappointed'
$ws.Range("D17").Value = '4.5'
$ws.Range("E17").Value = 'Very simple to use. Quite happy.'
$ws.Range("F17").Value = 'Easy to usestraightforward'
$ws.Range("G17").Value = 'for people who buy these leads  do not waste your time and money because if you are not membership you do not receive any help even if you pay for the lead 2 minutes after receiving the customer''s number he has changed his mind and does not want to do this job it is disgusting that Rated  People.  I do not return money for such cases I do not recommend it'
$ws.Range("H17").Value = 'Very simple to use'
